# Transition Rule summary tables update:
# Adds "Within 5 miles" and "Within 10 miles" of HFC production facility
# columns (F and G) to both the "Means" and "Standard Deviations" sheets,
# and updates the "Total Cancer Risk" / "Total Respiratory" rows (9 & 10)
# with refreshed values for columns B-E as well as the new F/G columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Means"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Means")

# New column headers
$ws.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws.Range("G1").Value = "Within 10 miles of HFC production facility"

# % White
$ws.Range("F2").Value = 92
$ws.Range("G2").Value = 91

# % Black or African American
$ws.Range("F3").Value = 1.7
$ws.Range("G3").Value = 2.3

# % Other
$ws.Range("F4").Value = 6.3
$ws.Range("G4").Value = 6.9

# % Hispanic
$ws.Range("F5").Value = 45
$ws.Range("G5").Value = 41

# Median Income [1,000 2019$]
$ws.Range("F6").Value = 69
$ws.Range("G6").Value = 59

# % Below Poverty Line
$ws.Range("F7").Value = 4.2
$ws.Range("G7").Value = 6.2

# % Below Half the Poverty Line
$ws.Range("F8").Value = 4.1
$ws.Range("G8").Value = 4.9

# Total Cancer Risk (per million) - values updated for B:E, new for F:G
$ws.Range("B9").Value = 26
$ws.Range("C9").Value = 28
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 20

# Total Respiratory (hazard quotient) - values updated for B:E, new for F:G
$ws.Range("B10").Value = 0.32
$ws.Range("C10").Value = 0.33
$ws.Range("D10").Value = 0.2
$ws.Range("E10").Value = 0.21
$ws.Range("F10").Value = 0.21
$ws.Range("G10").Value = 0.21

# ---------------------------------------------------------------------
# Sheet "Standard Deviations"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Standard Deviations")

# New column headers
$ws.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# % White
$ws.Range("F2").Value = 6.2
$ws.Range("G2").Value = 8.2

# % Black or African American
$ws.Range("F3").Value = 1.5
$ws.Range("G3").Value = 3

# % Other
$ws.Range("F4").Value = 5.9
$ws.Range("G4").Value = 8.1

# % Hispanic
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 23

# Median Income [1,000 2019$]
$ws.Range("F6").Value = 17
$ws.Range("G6").Value = 18

# % Below Poverty Line
$ws.Range("F7").Value = 5.2
$ws.Range("G7").Value = 8.5

# % Below Half the Poverty Line
$ws.Range("F8").Value = 3.6
$ws.Range("G8").Value = 8.6

# Total Cancer Risk (per million) - values updated for B:E, new for F:G
$ws.Range("B9").Value = 8.6
$ws.Range("C9").Value = 7.5
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# Total Respiratory (hazard quotient) - values updated for B:E, new for F:G
$ws.Range("B10").Value = 0.14
$ws.Range("C10").Value = 0.076
$ws.Range("D10").Value = 0.000000000000000028
$ws.Range("E10").Value = 0.032
$ws.Range("F10").Value = 0.029
$ws.Range("G10").Value = 0.038
